# Automatic update of files.
# - Bumps the "Förändrad" (Changed) date in column C for every data row.
# - Rewrites the municipality-based folder names used in the hyperlink
#   formulas (columns S-Y) to their numeric "Logging_<id>" equivalents
#   for the rows that moved to a new folder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the extent of the data so the script keeps working even if the
# sheet grows or shrinks a little.
$lastRow = $ws.Cells.Item(1, 1).SpecialCells(11).Row  # xlCellTypeLastCell -> Row

$firstDataRow = 2

# 1) Update column C (Förändrad) from 45207 to 45208 for every data row.
#    (Use Value2 - the COM interop layer's .Value accessor does not
#    reliably round-trip simple numerics in this environment.)
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# 2) Rewrite the folder name inside the HYPERLINK() formulas for the
#    rows whose logging case moved to a different (numeric) folder.
$folderRenames = @{
    2 = @{ Old = "Logging_ANGE"; New = "Logging_2260" }
    3 = @{ Old = "Logging_ANGE"; New = "Logging_2260" }
    4 = @{ Old = "Logging_ANGE"; New = "Logging_2260" }
    5 = @{ Old = "Logging_SOLLEFTEA"; New = "Logging_2283" }
}

$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

foreach ($row in $folderRenames.Keys) {
    $oldName = $folderRenames[$row].Old
    $newName = $folderRenames[$row].New
    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$row")
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            if ($formula.Contains($oldName)) {
                $cell.Formula = $formula.Replace($oldName, $newName)
            }
        }
    }
}
